$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (the numbering column). This shifts the question column
# (old B) into A and the answer column (old C) into B, carrying over all
# values, shared-string references and cell styles automatically.
$ws.Columns.Item(1).Delete()

# Mirror the resulting selection state: the whole (new) column A ends up
# selected, as happens right after deleting a column via its header.
$ws.Range("A1:A1048576").Select()
